# This script updates column G ("K") on Sheet1 of the active workbook.
# The underlying save_data generator was changed to compute "K" (strikeouts)
# differently (regenerated std/mean, calc and write s_vals) instead of the
# previous "Strike#" value, so only the G column values for the 28 data
# rows (rows 2-29) need to be rewritten with their newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G (K), keyed by row number (row 2 = first data row).
$newK = @{
    2  = 4
    3  = 8
    4  = 9
    5  = 2
    6  = 4
    7  = 8
    8  = 8
    9  = 6
    10 = 4
    11 = 4
    12 = 4
    13 = 5
    14 = 4
    15 = 8
    16 = 5
    17 = 7
    18 = 3
    19 = 6
    20 = 3
    21 = 4
    22 = 3
    23 = 5
    24 = 3
    25 = 3
    26 = 7
    27 = 2
    28 = 3
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
